$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append after the existing last row (105)
$data = @(
    @{ Row = 106; Date = 45901; B = 0.0671749355181308; C = 0.140943592796547 },
    @{ Row = 107; Date = 45931; B = 0.116275982400243;  C = 0.104947589990007 },
    @{ Row = 108; Date = 45962; B = 0.153997875891367;  C = 0.256181829779091 }
)

foreach ($item in $data) {
    $row = $item.Row

    # Copy the style/format of the last populated date cell (A105) so the new
    # date cell shares the same cell style (numFmtId 14 / s="1") rather than
    # creating a brand-new style entry.
    $ws.Range("A105").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $item.Date
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C
}

$excel.CutCopyMode = 0
